$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.752878
$ws.Range("N2").Value = 5.258634
$ws.Range("O2").Value = 0.1377607590022273
$ws.Range("P2").Value = 0.1377607590022273
$ws.Range("Q2").Value = 71.43572309359067
$ws.Range("R2").Value = 642.921507842316
$ws.Range("S2").Value = 0.002929533003583409
$ws.Range("T2").Value = 0.002929533003583409

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.076282333333333
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.2417686736584878
$ws.Range("P3").Value = 0.2417686736584878
$ws.Range("Q3").Value = 125.3689377821531
$ws.Range("R3").Value = 1128.320440039378
$ws.Range("S3").Value = 0.005141299408082351
$ws.Range("T3").Value = 0.005141299408082352

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.303088666666667
$ws.Range("N4").Value = 21.909266
$ws.Range("O4").Value = 0.5739583917309499
$ws.Range("P4").Value = 0.5739583917309499
$ws.Range("Q4").Value = 297.6256303746982
$ws.Range("R4").Value = 2678.630673372284
$ws.Range("S4").Value = 0.0122054354479296
$ws.Range("T4").Value = 0.0122054354479296

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5918243333333334
$ws.Range("N5").Value = 1.775473
$ws.Range("O5").Value = 0.04651217560833507
$ws.Range("P5").Value = 0.04651217560833507
$ws.Range("Q5").Value = 24.11884865692245
$ws.Range("R5").Value = 217.069637912302
$ws.Range("S5").Value = 0.0009890984522732036
$ws.Range("T5").Value = 0.0009890984522732036

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.752878
$ws.Range("N6").Value = 5.258634
$ws.Range("O6").Value = 0.1377607590022273
$ws.Range("P6").Value = 0.1377607590022273
$ws.Range("Q6").Value = 2961.118061291253
$ws.Range("R6").Value = 26650.06255162128
$ws.Range("S6").Value = 0.1214335449043415
$ws.Range("T6").Value = 0.1214335449043415

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.076282333333333
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.2417686736584878
$ws.Range("P7").Value = 0.2417686736584878
$ws.Range("Q7").Value = 5196.730849987582
$ws.Range("R7").Value = 46770.57764988824
$ws.Range("S7").Value = 0.2131145857631082
$ws.Range("T7").Value = 0.2131145857631082

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.303088666666667
$ws.Range("N8").Value = 21.909266
$ws.Range("O8").Value = 0.5739583917309499
$ws.Range("P8").Value = 0.5739583917309499
$ws.Range("Q8").Value = 12337.02959023852
$ws.Range("R8").Value = 111033.2663121467
$ws.Range("S8").Value = 0.5059336391603144
$ws.Range("T8").Value = 0.5059336391603144

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5918243333333334
$ws.Range("N9").Value = 1.775473
$ws.Range("O9").Value = 0.04651217560833507
$ws.Range("P9").Value = 0.04651217560833507
$ws.Range("Q9").Value = 999.762517725129
$ws.Range("R9").Value = 8997.862659526159
$ws.Range("S9").Value = 0.04099961706251963
$ws.Range("T9").Value = 0.04099961706251963

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.752878
$ws.Range("N10").Value = 5.258634
$ws.Range("O10").Value = 0.1377607590022273
$ws.Range("P10").Value = 0.1377607590022273
$ws.Range("Q10").Value = 176.9800325112487
$ws.Range("R10").Value = 1592.820292601238
$ws.Range("S10").Value = 0.007257837168354854
$ws.Range("T10").Value = 0.007257837168354854

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.076282333333333
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.2417686736584878
$ws.Range("P11").Value = 0.2417686736584878
$ws.Range("Q11").Value = 310.5980834759254
$ws.Range("R11").Value = 2795.382751283329
$ws.Range("S11").Value = 0.01273742739609948
$ws.Range("T11").Value = 0.01273742739609948

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.303088666666667
$ws.Range("N12").Value = 21.909266
$ws.Range("O12").Value = 0.5739583917309499
$ws.Range("P12").Value = 0.5739583917309499
$ws.Range("Q12").Value = 737.3592855060068
$ws.Range("R12").Value = 6636.233569554061
$ws.Range("S12").Value = 0.03023862948175767
$ws.Range("T12").Value = 0.03023862948175767

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.5918243333333334
$ws.Range("N13").Value = 1.775473
$ws.Range("O13").Value = 0.04651217560833507
$ws.Range("P13").Value = 0.04651217560833507
$ws.Range("Q13").Value = 59.75378192565678
$ws.Range("R13").Value = 537.784037330911
$ws.Range("S13").Value = 0.002450464118782653
$ws.Range("T13").Value = 0.002450464118782653

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.752878
$ws.Range("N14").Value = 5.258634
$ws.Range("O14").Value = 0.1377607590022273
$ws.Range("P14").Value = 0.1377607590022273
$ws.Range("Q14").Value = 149.7181257201587
$ws.Range("R14").Value = 1347.463131481428
$ws.Range("S14").Value = 0.006139843925947564
$ws.Range("T14").Value = 0.006139843925947564

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.076282333333333
$ws.Range("N15").Value = 9.228847
$ws.Range("O15").Value = 0.2417686736584878
$ws.Range("P15").Value = 0.2417686736584878
$ws.Range("Q15").Value = 262.7537256629971
$ws.Range("R15").Value = 2364.783530966974
$ws.Range("S15").Value = 0.01077536109119771
$ws.Range("T15").Value = 0.01077536109119771

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.303088666666667
$ws.Range("N16").Value = 21.909266
$ws.Range("O16").Value = 0.5739583917309499
$ws.Range("P16").Value = 0.5739583917309499
$ws.Range("Q16").Value = 623.7768670389303
$ws.Range("R16").Value = 5613.991803350372
$ws.Range("S16").Value = 0.0255806876409481
$ws.Range("T16").Value = 0.0255806876409481

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.5918243333333334
$ws.Range("N17").Value = 1.775473
$ws.Range("O17").Value = 0.04651217560833507
$ws.Range("P17").Value = 0.04651217560833507
$ws.Range("Q17").Value = 50.54934224871845
$ws.Range("R17").Value = 454.944080238466
$ws.Range("S17").Value = 0.002072995974759586
$ws.Range("T17").Value = 0.002072995974759586
